$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.745.03'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.272.06'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.66%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.35'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.69%  '

$ws.Range("E6").Value = '  +1.58%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '78.57'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +7.85%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.644'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.02'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.97%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0966'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.39'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.32%  '

$ws.Range("E13").Value = '  -0.62%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.611.95'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.08'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.867'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.36%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.279.71'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.16%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.648.61'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.68%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0994'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.36%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.21'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.01'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.39%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '233.39'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.17'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.84%  '

$ws.Range("E24").Value = '  -3.01%  '

$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.37'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.91%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.33'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.84%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.17'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.07'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.72'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.83'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.81%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0850'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.88%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.122'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.79%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.71'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.09%  '

$ws.Range("E35").Value = '  +0.09%  '

$ws.Range("E36").Value = '  -5.52%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.76'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.76%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0302'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.87%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.48'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.46%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.26'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.60%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.96'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.45%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '115.12'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +17.64%  '

$ws.Range("E43").Value = '  -2.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '61.32'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.54%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.89'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.50%  '

$ws.Range("E46").Value = '  -2.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.63'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.76%  '

$ws.Range("E48").Value = '  -0.19%  '

$ws.Range("E49").Value = '  -4.36%  '

$ws.Range("E50").Value = '  -2.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.23'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.06%  '
